# Append 4 new daily rows (230-233) to the report, covering the next
# four dates after the last existing row (44303 -> 44304..44307),
# keeping the same layout/formatting as the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (custom date/time number format,
# borders, centered alignment, etc.) from the last existing row (A229)
# down onto the new date cells, so the new rows look like the rest of
# the table.
$ws.Range("A229").Copy() | Out-Null
$ws.Range("A230:A233").PasteSpecial(-4122) | Out-Null

# New data to append: date serial, nuovi pos., somma mobile 7gg.,
# somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44304, 1,  33, 304.0073698756333),
    @(44305, 11, 41, 377.7061262091202),
    @(44306, 3,  38, 350.0690925840626),
    @(44307, 0,  37, 340.8567480423768)
)

$r = 230
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
